$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CAPEX input value for BESS (B2): 238 -> 275
$ws.Range("B2").Value = 275

# Move the active selection to C2 (matches author's final cursor position)
$ws.Range("C2").Select()
